$wb = $excel.ActiveWorkbook

# "Raw Data" sheet: fill in the previously-empty row 38 with the new
# product line, and make sure an (empty) row 39 exists below it so the
# sheet's used range grows to A1:B39.
$wsRaw = $wb.Worksheets.Item("Raw Data")
$wsRaw.Range("A38").Value = "165349_B_Framside"
$wsRaw.Range("B38").Value = 4419.727274000001

# Touch row 39 (without actually changing anything visible) just so the
# worksheet registers it as part of the used range, matching the new
# dimension ref="A1:B39" / trailing empty <row r="39"/>.
$wsRaw.Cells.Item(39, 1).Font.Bold = $false

# "Results" sheet: the machining time for product 165349 was recalculated.
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Range("B18").Value = "1:48:17.677436000001762"
